$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1 / sheet1.xml) - 16 cell updates in column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 54
$ws1.Range("F4").Value = 1105
$ws1.Range("F5").Value = 369
$ws1.Range("F7").Value = 590
$ws1.Range("F8").Value = 1513
$ws1.Range("F11").Value = 3054
$ws1.Range("F12").Value = 565
$ws1.Range("F13").Value = 1725
$ws1.Range("F17").Value = 1449
$ws1.Range("F20").Value = 1179
$ws1.Range("F22").Value = 431
$ws1.Range("F23").Value = 60
$ws1.Range("F24").Value = 4662
$ws1.Range("F25").Value = 732
$ws1.Range("F28").Value = 43
$ws1.Range("F29").Value = 88

# Sheet "演出" (sheetId 2 / sheet2.xml) - 1 cell update
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value = 50

# Sheet "本地生活" (sheetId 3 / sheet3.xml) - 1 cell update
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 30

# Sheet "全部类型" (sheetId 4 / sheet4.xml) - 18 cell updates in column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 54
$ws4.Range("F4").Value = 30
$ws4.Range("F13").Value = 50
$ws4.Range("F15").Value = 1105
$ws4.Range("F16").Value = 369
$ws4.Range("F18").Value = 590
$ws4.Range("F19").Value = 1513
$ws4.Range("F22").Value = 3054
$ws4.Range("F23").Value = 565
$ws4.Range("F24").Value = 1725
$ws4.Range("F28").Value = 1449
$ws4.Range("F33").Value = 1179
$ws4.Range("F35").Value = 431
$ws4.Range("F36").Value = 60
$ws4.Range("F37").Value = 4662
$ws4.Range("F38").Value = 732
$ws4.Range("F43").Value = 43
$ws4.Range("F44").Value = 88
